$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly refresh of the "Coco" price series: every existing observation in
# rows 33-42 takes on the values previously held by the row above it (the
# oldest row 32 observation rolls down into the brand-new row 43), and row
# 32 is populated with the newest observation.

$ws.Range("D32").Value = 44466
$ws.Range("M32").Value = 70

$ws.Range("D33").Value = 44389
$ws.Range("M33").Value = 20
$ws.Range("N33").Value = 24000
$ws.Range("O33").Value = 24000
$ws.Range("P33").Value = 24000
$ws.Range("S33").Value = 1200

$ws.Range("D34").Value = 44249
$ws.Range("M34").Value = 15
$ws.Range("N34").Value = 25000
$ws.Range("O34").Value = 25000
$ws.Range("P34").Value = 25000
$ws.Range("S34").Value = 1250

$ws.Range("D35").Value = 44390
$ws.Range("M35").Value = 10
$ws.Range("N35").Value = 24000
$ws.Range("O35").Value = 24000
$ws.Range("P35").Value = 24000
$ws.Range("S35").Value = 1200

$ws.Range("D36").Value = 44251
$ws.Range("M36").Value = 15
$ws.Range("N36").Value = 25000
$ws.Range("O36").Value = 25000
$ws.Range("P36").Value = 25000
$ws.Range("S36").Value = 1250

$ws.Range("D37").Value = 44433
$ws.Range("M37").Value = 10
$ws.Range("N37").Value = 24000
$ws.Range("O37").Value = 24000
$ws.Range("P37").Value = 24000
$ws.Range("S37").Value = 1200

$ws.Range("D38").Value = 44221
$ws.Range("M38").Value = 30
$ws.Range("N38").Value = 25000
$ws.Range("O38").Value = 25000
$ws.Range("P38").Value = 25000
$ws.Range("S38").Value = 1250

$ws.Range("D39").Value = 44363
$ws.Range("M39").Value = 30
$ws.Range("N39").Value = 24000
$ws.Range("O39").Value = 24000
$ws.Range("P39").Value = 24000
$ws.Range("S39").Value = 1200

$ws.Range("D40").Value = 44356
$ws.Range("M40").Value = 15
$ws.Range("N40").Value = 24000
$ws.Range("O40").Value = 24000
$ws.Range("P40").Value = 24000
$ws.Range("S40").Value = 1200

$ws.Range("D41").Value = 44175
$ws.Range("M41").Value = 25
$ws.Range("N41").Value = 23000
$ws.Range("O41").Value = 23000
$ws.Range("P41").Value = 23000
$ws.Range("S41").Value = 1150

$ws.Range("D42").Value = 44461
$ws.Range("M42").Value = 30
$ws.Range("N42").Value = 24000
$ws.Range("O42").Value = 24000
$ws.Range("P42").Value = 24000
$ws.Range("S42").Value = 1200

# New row 43, carrying the values that used to be in row 42.
$ws.Range("A43").Value = 10
$ws.Range("B43").Value = "Vega Modelo de Temuco"
$ws.Range("C43").Value = "La Araucanía"
$ws.Range("D43").Value = 44425
$ws.Range("D43").NumberFormat = $ws.Range("D42").NumberFormat
$ws.Range("E43").Value = 9
$ws.Range("F43").Value = "Fruta"
$ws.Range("G43").Value = 100108
$ws.Range("H43").Value = "Tropicales y subtropicales"
$ws.Range("I43").Value = 100108007
$ws.Range("J43").Value = "Coco"
$ws.Range("K43").Value = "Sin especificar"
$ws.Range("L43").Value = "Primera"
$ws.Range("M43").Value = 15
$ws.Range("N43").Value = 24000
$ws.Range("O43").Value = 24000
$ws.Range("P43").Value = 24000
$ws.Range("Q43").Value = "$/malla 20 unidades"
$ws.Range("R43").Value = "Perú"
$ws.Range("S43").Value = 1200
$ws.Range("T43").Value = 20
